# Update the two "備註說明" (remarks) cells on the DBD sheet so the
# code/label order in the note text is swapped (code first, then label).
$wb = $excel.ActiveWorkbook

$wsDBD = $wb.Worksheets.Item("DBD")

# G11: AcSubBookCode note — was "一般：00A\n利變：201" now "00A:一般\n201:利變"
$wsDBD.Range("G11").Value = "00A:一般`n201:利變"

# G10: AssetClassNo note — was "一類：11、12\n..." now "11/12:一類\n..."
$wsDBD.Range("G10").Value = "11/12:一類`n21/22/23:二類`n3:三類`n4:四類`n5:五類`n6:折溢價與催收`n7:應收利息提列"

# The active sheet moves back to "DBD" (first tab) and the DBD sheet view
# scrolls down a bit with the selection now on G11 instead of G10.
$wsDBD.Activate()
$wsDBD.Range("G11").Select()
